$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 2-8 (Orange/shutterstock_513244762 comparisons).
# The data that used to live in rows 9-11 (Dog-2/Dog-1/Cat-1 comparisons)
# shifts up to become the new rows 2-4, leaving a table of just
# A1:C4 (header + 3 data rows).
$ws.Rows("2:8").Delete()
